$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B2').Formula = 'Bitcoin'
$ws.Range('C2').Formula = 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc'
$ws.Range('D2').Formula = '42.358.49'
$ws.Range('E2').Formula = '  +1.46%  '

$ws.Range('B3').Formula = 'Ethereum'
$ws.Range('C3').Formula = 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth'
$ws.Range('D3').Formula = '2.289.63'
$ws.Range('E3').Formula = '  +0.40%  '

$ws.Range('B4').Formula = 'TetherUSD'
$ws.Range('C4').Formula = 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt'
$ws.Range('D4').Value = "'1.00"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Formula = '  +0.16%  '

$ws.Range('B5').Formula = 'BNB'
$ws.Range('C5').Formula = 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'
$ws.Range('D5').Value = "'314.79"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Formula = '  +1.34%  '

$ws.Range('B6').Formula = 'Solana'
$ws.Range('C6').Formula = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D6').Value = "'102.35"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Formula = '  -2.93%  '

$ws.Range('B7').Formula = 'XRP'
$ws.Range('C7').Formula = 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp'
$ws.Range('D7').Value = "'0.628"
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Formula = '  +0.46%  '

$ws.Range('B8').Formula = 'USDC'
$ws.Range('C8').Formula = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').Value = "'1.00"
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Formula = '  +0.23%  '

$ws.Range('B9').Formula = 'Cardano'
$ws.Range('C9').Formula = 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
$ws.Range('D9').Value = "'0.603"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Formula = '  -0.46%  '

$ws.Range('B10').Formula = 'Avalanche'
$ws.Range('C10').Formula = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D10').Value = "'39.41"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Formula = '  -2.16%  '

$ws.Range('B11').Formula = 'Dogecoin'
$ws.Range('C11').Formula = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').Value = "'0.0905"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Formula = '  -0.45%  '

$ws.Range('B12').Formula = 'Polkadot'
$ws.Range('C12').Formula = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range('D12').Value = "'8.34"
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Formula = '  +1.33%  '

$ws.Range('B13').Formula = 'TRON'
$ws.Range('C13').Formula = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D13').Value = "'0.106"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Formula = '  +0.63%  '

$ws.Range('B14').Formula = 'Polygon'
$ws.Range('C14').Formula = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D14').Value = "'0.956"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Formula = '  -0.82%  '

$ws.Range('B15').Formula = 'Chainlink'
$ws.Range('C15').Formula = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = "'15.16"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Formula = '  -1.64%  '

$ws.Range('B16').Formula = 'WrappedliquidstakedEther2.0'
$ws.Range('C16').Formula = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D16').Formula = '2.642.00'
$ws.Range('E16').Formula = '  +0.51%  '

$ws.Range('B17').Formula = 'WrappedEther'
$ws.Range('C17').Formula = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Formula = '2.290.76'
$ws.Range('E17').Formula = '  +0.07%  '

$ws.Range('B18').Formula = 'WrappedBTC'
$ws.Range('C18').Formula = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D18').Formula = '42.308.71'
$ws.Range('E18').Formula = '  +1.18%  '

$ws.Range('B19').Formula = 'Uniswap'
$ws.Range('C19').Formula = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D19').Value = "'7.36"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Formula = '  -1.59%  '

$ws.Range('B20').Formula = 'ShibaInu'
$ws.Range('C20').Formula = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = "'0.0000105"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Formula = '  +0.57%  '

$ws.Range('B21').Formula = 'Litecoin'
$ws.Range('C21').Formula = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D21').Value = "'73.04"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Formula = '  -0.19%  '

$ws.Range('B22').Formula = 'InternetComputer(DFINITY)'
$ws.Range('C22').Formula = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D22').Value = "'11.69"
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Formula = '  +25.84%  '

$ws.Range('B23').Formula = 'PancakeSwap'
$ws.Range('C23').Formula = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D23').Value = "'3.52"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Formula = '  +2.36%  '

$ws.Range('B24').Formula = 'BitcoinCash'
$ws.Range('C24').Formula = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').Value = "'274.51"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Formula = '  +7.44%  '

$ws.Range('B25').Formula = 'ImmutableX'
$ws.Range('C25').Formula = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D25').Value = "'2.25"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Formula = '  -2.33%  '

$ws.Range('B26').Formula = 'Dai'
$ws.Range('C26').Formula = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D26').Value = "'1.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Formula = '  -0.28%  '

$ws.Range('B27').Formula = 'Cosmos'
$ws.Range('C27').Formula = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range('D27').Value = "'10.77"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Formula = '  -1.41%  '

$ws.Range('B28').Formula = 'Toncoin'
$ws.Range('C28').Formula = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D28').Value = "'2.36"
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Formula = '  +3.41%  '

$ws.Range('B29').Formula = 'EthereumClassic'
$ws.Range('C29').Formula = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D29').Value = "'22.65"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Formula = '  +0.15%  '

$ws.Range('B30').Formula = 'InjectiveProtocol'
$ws.Range('C30').Formula = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D30').Value = "'37.34"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Formula = '  +5.27%  '

$ws.Range('B31').Formula = 'Monero'
$ws.Range('C31').Formula = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D31').Value = "'165.52"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Formula = '  -0.24%  '

$ws.Range('B32').Formula = 'Hedera'
$ws.Range('C32').Formula = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').Value = "'0.0868"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Formula = '  -2.26%  '

$ws.Range('B33').Formula = 'Filecoin'
$ws.Range('C33').Formula = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = "'5.92"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Formula = '  +3.27%  '

$ws.Range('B34').Formula = 'Stellar'
$ws.Range('C34').Formula = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D34').Value = "'0.133"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Formula = '  +2.81%  '

$ws.Range('B35').Formula = 'WEMIXToken'
$ws.Range('C35').Formula = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D35').Value = "'2.65"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Formula = '  -9.14%  '

$ws.Range('B36').Formula = 'Kaspa'
$ws.Range('C36').Formula = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D36').Value = "'0.118"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Formula = '  -0.13%  '

$ws.Range('B37').Formula = 'RenderToken'
$ws.Range('C37').Formula = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = "'4.54"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Formula = '  -0.40%  '

$ws.Range('B38').Formula = 'VeChain'
$ws.Range('C38').Formula = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = "'0.0362"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Formula = '  +3.05%  '

$ws.Range('B39').Formula = 'NEARProtocol'
$ws.Range('C39').Formula = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').Value = "'3.70"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Formula = '  +2.36%  '

$ws.Range('B40').Formula = 'LidoDAOToken'
$ws.Range('C40').Formula = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').Value = "'2.75"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Formula = '  -1.44%  '

$ws.Range('B41').Formula = 'ARBITRUM'
$ws.Range('C41').Formula = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D41').Value = "'1.49"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Formula = '  +2.32%  '

$ws.Range('B42').Formula = 'BitcoinSV'
$ws.Range('C42').Formula = 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
$ws.Range('D42').Value = "'95.80"
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Formula = '  -2.22%  '

$ws.Range('B43').Formula = 'MultiversX'
$ws.Range('C43').Formula = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$ws.Range('D43').Value = "'69.54"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Formula = '  -2.11%  '

$ws.Range('B44').Formula = 'FirstDigitalUSD'
$ws.Range('C44').Formula = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').Value = "'1.00"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Formula = '  -0.09%  '

$ws.Range('B45').Formula = 'Algorand'
$ws.Range('C45').Formula = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').Value = "'0.224"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Formula = '  -0.89%  '

$ws.Range('B46').Formula = 'Celestia'
$ws.Range('C46').Formula = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D46').Value = "'11.96"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Formula = '  -2.30%  '

$ws.Range('B47').Formula = 'ordi'
$ws.Range('C47').Formula = 'https://coinranking.com/coin/j7-7vPrOi+ordi-ordi'
$ws.Range('D47').Value = "'79.91"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Formula = '  +6.09%  '

$ws.Range('B48').Formula = 'Aave'
$ws.Range('C48').Formula = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = "'112.27"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Formula = '  +0.16%  '

$ws.Range('B49').Formula = 'FraxShare'
$ws.Range('C49').Formula = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D49').Value = "'8.95"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Formula = '  -0.96%  '

$ws.Range('B50').Formula = 'THORChain'
$ws.Range('C50').Formula = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D50').Value = "'5.25"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Formula = '  -0.48%  '

$ws.Range('B51').Formula = 'Maker'
$ws.Range('C51').Formula = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D51').Formula = '1.586.20'
$ws.Range('E51').Formula = '  +2.00%  '
